$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 16-28 on columns C (doc number), D (name), E (period), F (valor mora).
# New order: Randy / Maira / Jorge-2502 moved to the top, followed by Jorge's
# remaining periods in descending order (2501 .. 2404).
$rows = @(
    @{ Row = 16; C = "1128050520"; D = "RANDY JAVIER TORRENTE HANNA";     E = "2502"; F = 32933 },
    @{ Row = 17; C = "33102541";   D = "MAIRA ALEJANDRA MARTINEZ ZUÑIGA"; E = "2502"; F = 32933 },
    @{ Row = 18; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2502"; F = 32933 },
    @{ Row = 19; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2501"; F = 52000 },
    @{ Row = 20; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2412"; F = 52000 },
    @{ Row = 21; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2411"; F = 52000 },
    @{ Row = 22; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2410"; F = 52000 },
    @{ Row = 23; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2409"; F = 52000 },
    @{ Row = 24; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2408"; F = 52000 },
    @{ Row = 25; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2407"; F = 52000 },
    @{ Row = 26; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2406"; F = 52000 },
    @{ Row = 27; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2405"; F = 52000 },
    @{ Row = 28; C = "1041977150"; D = "JORGE HUMBERTO RAMIREZ MARTINEZ"; E = "2404"; F = 52000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
